# edit.ps1 - applies the "Fixed test specs and added last documentation" edit
#   1) Date: 25/10/2019  ->  Date: 25/09/2019   (split into extra runs exactly
#      as produced by a user typing "/" before the _GoBack bookmark and then
#      fixing "10" to "09" just after it)
#   2) Merge the three runs that make up the "Search for a medical code..."
#      bullet into a single run (no functional text change, just a retype).

$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Change 1: fix the date from 25/10/2019 to 25/09/2019
# -------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")

# The bookmark currently sits directly in front of the run containing
# "/10/2019". Replace that run's text in place first (keeps a single run),
# then split "09" away from "/2019", then insert the new "/" run in front
# of the bookmark. A harmless Bold on/off toggle is used after each new
# piece of text so the engine keeps it as its own <w:r> instead of folding
# it back into the neighbouring run that has identical formatting.

$oldRange = $d.Range($bm.End, $bm.End + 8)
if ($oldRange.Text -eq "/10/2019") {
    $oldRange.Text = "09/2019"
}

$bm = $d.Bookmarks("_GoBack")
$slashYear = $d.Range($bm.End + 2, $bm.End + 7)
if ($slashYear.Text -eq "/2019") {
    $slashYear.Font.Bold = 1
    $slashYear.Font.Bold = 0
}

$bm = $d.Bookmarks("_GoBack")
$insertionPoint = $d.Range($bm.Start, $bm.Start)
$insertionPoint.InsertBefore("/")

$bm = $d.Bookmarks("_GoBack")
$newSlash = $d.Range($bm.Start - 1, $bm.Start)
if ($newSlash.Text -eq "/") {
    $newSlash.Font.Bold = 1
    $newSlash.Font.Bold = 0
}

# -------------------------------------------------------------------------
# Change 2: merge the three runs of the "Search for a medical code..."
# bullet into a single run with the same combined text.
# -------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Search for a medical code in the procedure search bar (Eg. 065 & 039 are known working codes).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Search for a medical code in the procedure search bar (Eg. 065 & 039 are known working codes).",
    2
)
